# Weekly price-sheet update: a new daily observation row is inserted
# right before the existing row 172 (pushing the rest of the table down
# by one row), matching the "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 172..189 down to 173..190 and open up a fresh row 172.
$ws.Rows(172).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(172, 1).Value  = 7
$ws.Cells.Item(172, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(172, 3).Value  = 'Ñuble'
$ws.Cells.Item(172, 4).Value  = 45142
$ws.Cells.Item(172, 5).Value  = 16
$ws.Cells.Item(172, 6).Value  = 100112037
$ws.Cells.Item(172, 7).Value  = 'Cebollín'
$ws.Cells.Item(172, 8).Value  = 'Sin especificar'
$ws.Cells.Item(172, 9).Value  = 'Primera'
$ws.Cells.Item(172, 10).Value = 60
$ws.Cells.Item(172, 11).Value = 6000
$ws.Cells.Item(172, 12).Value = 6000
$ws.Cells.Item(172, 13).Value = 6000
$ws.Cells.Item(172, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(172, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(172, 16).Value = 167
$ws.Cells.Item(172, 17).Value = 36
$ws.Cells.Item(172, 18).Value = 'Hortaliza'
